$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was inserted as row 12, pushing the previously
# existing rows 12..103 down to 13..104 (dimension grows from R103 to R104).
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new observation's data.
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = 44831
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = 100112026
$ws.Cells.Item(12, 7).Value = "Haba"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 80
$ws.Cells.Item(12, 11).Value = 13500
$ws.Cells.Item(12, 12).Value = 13500
$ws.Cells.Item(12, 13).Value = 13500
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 540
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
